$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column H ("Absent") had been left blank/inconsistent for a few rows
# while consolidating the report. Fill in the computed Absent flag for
# the remaining rows (10-13) so the report is fully consolidated.
$ws.Range("H10").Value = 1
$ws.Range("H11").Value = 0
$ws.Range("H12").Value = 1
$ws.Range("H13").Value = 0
